$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 63, pushing existing rows 63-117 down to 65-119.
$ws.Rows.Item(63).Insert()
$ws.Rows.Item(63).Insert()

# New row 63 data
$ws.Cells.Item(63, 1).Value = 5
$ws.Cells.Item(63, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(63, 3).Value = "Maule"
$ws.Cells.Item(63, 4).Value = 44634
$ws.Cells.Item(63, 5).Value = 7
$ws.Cells.Item(63, 6).Value = "Fruta"
$ws.Cells.Item(63, 7).Value = 100108
$ws.Cells.Item(63, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(63, 9).Value = 100108002
$ws.Cells.Item(63, 10).Value = "Mango"
$ws.Cells.Item(63, 11).Value = "Sin especificar"
$ws.Cells.Item(63, 12).Value = "Especial"
$ws.Cells.Item(63, 13).Value = 210
$ws.Cells.Item(63, 14).Value = 7000
$ws.Cells.Item(63, 15).Value = 7000
$ws.Cells.Item(63, 16).Value = 7000
$ws.Cells.Item(63, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(63, 18).Value = "Perú"
$ws.Cells.Item(63, 19).Value = 1750
$ws.Cells.Item(63, 20).Value = 4

# New row 64 data
$ws.Cells.Item(64, 1).Value = 5
$ws.Cells.Item(64, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(64, 3).Value = "Maule"
$ws.Cells.Item(64, 4).Value = 44634
$ws.Cells.Item(64, 5).Value = 7
$ws.Cells.Item(64, 6).Value = "Fruta"
$ws.Cells.Item(64, 7).Value = 100108
$ws.Cells.Item(64, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(64, 9).Value = 100108002
$ws.Cells.Item(64, 10).Value = "Mango"
$ws.Cells.Item(64, 11).Value = "Sin especificar"
$ws.Cells.Item(64, 12).Value = "Primera"
$ws.Cells.Item(64, 13).Value = 50
$ws.Cells.Item(64, 14).Value = 7000
$ws.Cells.Item(64, 15).Value = 7000
$ws.Cells.Item(64, 16).Value = 7000
$ws.Cells.Item(64, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(64, 18).Value = "Perú"
$ws.Cells.Item(64, 19).Value = 1750
$ws.Cells.Item(64, 20).Value = 4
